$d = $word.ActiveDocument
$vt = [char]11

# Para 6: Objetivos (PT) text -> "Programa resumido" (PT) short text
$d.Paragraphs.Item(6).Range.Text = 'A. Microeconomia. B. Macroeconomia. C. Desenvolvimento Econômico. D. Economia Internacional. E. Economia Brasileira'

# Para 7: Objetivos (EN, italic) text -> "Programa resumido" (EN) short text
$d.Paragraphs.Item(7).Range.Text = 'A. Microeconomy. B. Macroeconomy. C. Economic Development. D. International Economy. E. Brazilian Economy'

# Para 9: Docente name -> Objetivos (PT) paragraph text (ListBullet style kept)
$d.Paragraphs.Item(9).Range.Text = 'Apresentar os conceitos básicos da Ciência Econômica, capacitando os alunos a compreender os principais conceitos micro e macroeconômicos, e a analisar o discurso e a prática da economia, orientados pelo seu próprio senso crítico.'

# Para 11: "Programa resumido" (PT) short text -> full PT program (5 segments, line breaks)
$d.Paragraphs.Item(11).Range.Text = 'A. MICROECONOMIA: 1. Introdução aos conceitos de Economia e fundamentos da análise microeconômica. 2. Teoria do consumidor e da demanda. 3. Teoria da firma e da oferta. 4. Custos e formação de preços. 5. Estruturas de Mercado 6. Comportamento estratégico e concorrência. 7. Tecnologia como fator de produção. 8. Sustentabilidade: recursos, custos e indicadores ambientais. ' + $vt + 'B. MACROECONOMIA: 1. Fundamentos da análise macroeconômica. 2. Contabilidade nacional. 3. Equilíbrios clássicos e keynesiano. 4. Sistema monetário. 5. Política fiscal. 6. Economia mundial e comércio internacional. 7. Fundamentos da regressão como ferramenta para quantificar relações econômicas. 8. Setor público. ' + $vt + 'C. DESENVOLVIMENTO ECONÔMICO: 1. Fatores de Crescimento. 2. Fontes de Desenvolvimento. 3. Financiamento do Desenvolvimento Econômico. 4. Um modelo de Crescimento Econômico. 5. O Processo de internacionalização e globalização.' + $vt + 'D. ECONOMIA INTERNACIONAL: 1. Fundamentos do Comércio Internacional. 2. Determinação das Taxas de Câmbio. 3. Políticas Externas. 4. Fatores determinantes do comportamento das importações e exportações.' + $vt + 'E. ECONOMIA BRASILEIRA: 1. A experiência histórica da industrialização brasileira. 2. A internacionalização da economia brasileira. 3. Teoria dos ciclos e realidade brasileira. 4. Os ciclos econômicos do Brasil ao longo de sua história recente.'

# Para 12: "Programa resumido" (EN, italic) short text -> Objetivos (EN) text
$d.Paragraphs.Item(12).Range.Text = 'Introduce the basic concepts of Economic Science, enabling the students to understand the main micro and macroeconomic concepts and to analyze the discourse and practice of economics, guided by their critical sense.'

# Para 14: full PT program text -> old "Metodo" value text
$d.Paragraphs.Item(14).Range.Text = 'Provas, trabalhos em grupo, exercícios individuais e seminários.'

# Para 17: Avaliacao list - shift Metodo/Criterio/Norma values (process in reverse to avoid collisions)
$d.Paragraphs.Item(17).Range.Find.Execute('NF = (MF + PR)/2, onde MF é a média final da avaliação e PR é uma prova de recuperação.', $true, $false, $false, $false, $false, $true, 1, $false, 'MANKIW, N.G. Introdução à economia. São Paulo: Thomson Learning, 2006.' + $vt + $vt + 'SAMUELSON, P. Introdução à Economia. New York: Mc Graw-Hill Book Company.' + $vt + $vt + 'BACHA, Edmar. Introdução à Macroeconomia: Uma perspectiva brasileira. Rio de Janeiro: Campus, 1987.' + $vt + $vt + 'BACHA et al. Estado da Economia Mundial - Desafios e Respostas - Seminário em Homenagem a Pedro Malan. São Paulo: LTC, 2015.' + $vt + $vt + 'FURTADO, C. Formação econômica do Brasil. São Paulo: Companhia Editora Nacional, 2003.' + $vt + $vt + 'GREMAUD, A. P.; VASCONCELLOS, M. A. S.; TONETO JÚNIOR, R. Economia Brasileira Contemporânea. 8 ed. São Paulo: Atlas, 2017.' + $vt + $vt + 'VASCONCELLOS, M. A. S.; GARCIA, M. E. Fundamentos de Economia. 6 ed. São Paulo: Saraiva, 2018.' + $vt + $vt + 'VASCONCELLOS, M. A. S. ECONOMIA: Micro e Macro. São Paulo: Atlas, 2015.', 2) | Out-Null
$d.Paragraphs.Item(17).Range.Find.Execute('Média das atividades avaliativas.', $true, $false, $false, $false, $false, $true, 1, $false, 'NF = (MF + PR)/2, onde MF é a média final da avaliação e PR é uma prova de recuperação.', 2) | Out-Null
$d.Paragraphs.Item(17).Range.Find.Execute('Provas, trabalhos em grupo, exercícios individuais e seminários.', $true, $false, $false, $false, $false, $true, 1, $false, 'Média das atividades avaliativas.', 2) | Out-Null

# Para 19: bibliography text -> Docente name
$d.Paragraphs.Item(19).Range.Text = '7811306 - Diogo Ferraz'

Write-Host "Edit complete"
